# Speaker-notes formatting cleanup for slide 1:
#  - "Speech" + "/Story: " runs become one run's worth of text "Speech/Story: "
#  - drop the blank paragraph right after the "Have an actual story..." paragraph
#  - drop the blank paragraph + the "Interface: ..." paragraph that used to
#    follow "Multiple sections: ..." (the single trailing blank paragraph stays)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
$notesShape = $np.Shapes.Item(2)
$tr = $notesShape.TextFrame.TextRange

# PowerPoint exposes paragraph breaks in TextRange.Text as CR ("`r").
$paragraphs = $tr.Text -split "`r"

$result = @()

for ($i = 0; $i -lt $paragraphs.Length; $i++) {
    $text = $paragraphs[$i]

    if ($text.Trim().Length -eq 0 -and $i -gt 0) {
        $prev = $paragraphs[$i - 1]
        $next = if ($i + 1 -lt $paragraphs.Length) { $paragraphs[$i + 1] } else { "" }

        if ($prev -like "Have an actual story*") {
            # drop the blank paragraph after the story paragraph
            continue
        }
        if ($prev -like "Multiple sections*" -and $next -like "Interface:*") {
            # drop the blank paragraph that leads into the removed Interface paragraph
            continue
        }
    }

    if ($text -like "Interface: line on dataset*") {
        # this whole paragraph was removed
        continue
    }

    $result += $text
}

$tr.Text = [string]::Join("`n", $result)
